# Auto-generated from diff: update cryptos price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.915.94"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.224.71"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'241.94"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'73.00"
$ws.Range("E7").Value = "  -5.62%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = "  -4.74%  "
$ws.Range("D10").Value = "'42.49"
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("D11").Value = "'0.0953"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'6.98"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "2.559.31"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'14.31"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").Value = "'0.838"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("D17").Value = "2.205.39"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").Value = "41.755.46"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'0.0000106"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.19"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'72.70"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "'11.28"
$ws.Range("E22").Value = "  +23.10%  "
$ws.Range("D23").Value = "'229.83"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("E24").Value = "  -7.98%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'11.42"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "'166.99"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "'20.50"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'0.0799"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").Value = "'5.55"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").Value = "'30.34"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E36").Value = "  -10.00%  "
$ws.Range("E37").Value = "  -6.31%  "
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "'13.39"
$ws.Range("E39").Value = "  -7.72%  "
$ws.Range("D40").Value = "'2.13"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "'65.24"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "'5.65"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").Value = "'0.198"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").Value = "'8.74"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "'103.86"
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("D46").Value = "'0.101"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "'2.34"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").Value = "'1.17"
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").Value = "'2.69"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "2.427.82"
$ws.Range("E51").Value = "  -1.61%  "
